$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.263.30"
$ws.Range("E2").Value = "  -6.42%  "
$ws.Range("D3").Value = "2.910.72"
$ws.Range("E3").Value = "  -10.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "469.63"
$ws.Range("E5").Value = "  -12.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.45"
$ws.Range("E6").Value = "  -9.75%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "2.908.08"
$ws.Range("E8").Value = "  -10.15%  "
$ws.Range("E9").Value = "  -12.73%  "
$ws.Range("E10").Value = "  -14.34%  "
$ws.Range("E11").Value = "  -17.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  -17.80%  "
$ws.Range("E13").Value = "  -3.75%  "
$ws.Range("D14").Value = "3.414.62"
$ws.Range("E14").Value = "  -9.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.32"
$ws.Range("E15").Value = "  -14.80%  "
$ws.Range("D16").Value = "55.269.96"
$ws.Range("E16").Value = "  -6.48%  "
$ws.Range("D17").Value = "2.921.74"
$ws.Range("E17").Value = "  -9.66%  "
$ws.Range("E18").Value = "  -17.52%  "
$ws.Range("E19").Value = "  -14.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.38"
$ws.Range("E20").Value = "  -14.02%  "
$ws.Range("E21").Value = "  -15.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "304.66"
$ws.Range("E22").Value = "  -15.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  -15.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "58.74"
$ws.Range("E25").Value = "  -16.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.153"
$ws.Range("E27").Value = "  -9.92%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -19.23%  "
$ws.Range("E30").Value = "  -17.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("E31").Value = "  -10.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.85"
$ws.Range("E32").Value = "  -14.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.95"
$ws.Range("E33").Value = "  -15.58%  "
$ws.Range("E34").Value = "  -18.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "143.49"
$ws.Range("E35").Value = "  -12.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.14"
$ws.Range("E36").Value = "  -16.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.33"
$ws.Range("E37").Value = "  -16.76%  "
$ws.Range("E38").Value = "  -16.71%  "
$ws.Range("D39").Value = "2.941.86"
$ws.Range("E39").Value = "  -9.93%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0601"
$ws.Range("E41").Value = "  -15.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.21"
$ws.Range("E42").Value = "  -18.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.92"
$ws.Range("E43").Value = "  -15.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.945"
$ws.Range("E44").Value = "  -14.01%  "
$ws.Range("E45").Value = "  -16.82%  "
$ws.Range("E46").Value = "  -16.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.28"
$ws.Range("E47").Value = "  -15.07%  "
$ws.Range("D48").Value = "2.021.70"
$ws.Range("E48").Value = "  -12.01%  "
$ws.Range("E49").Value = "  -16.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.29"
$ws.Range("E50").Value = "  -16.62%  "
$ws.Range("E51").Value = "  -14.72%  "
